# Nexial xml-showcase.xlsx: add a new "localdb" command-type category to the
# hidden '#system' lookup sheet.
#
# This mirrors what happens in the real workbook when a new column is
# inserted immediately before the existing "macro" column (column N) on the
# '#system' sheet, the new column is populated with the "localdb" category
# header plus its six command names, and the category name "localdb" is
# inserted (in alphabetical order) into the "target" list in column A.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("#system")

# 1) Insert a new column at N. Everything that used to live in columns
#    N..AC (macro, mail, number, pdf, rdbms, redis, sms, sound, ssh, step,
#    web, webalert, webcookie, ws, ws.async, xml) shifts right by one column,
#    into O..AD, leaving a blank column N for the new "localdb" data.
$ws.Columns("N").Insert()

# 2) Populate the new column N with the "localdb" header and its commands.
$localdb = @(
    "localdb",
    "cloneTable(var,source,target)",
    "dropTables(var,tables)",
    "exportCSV(sql,output)",
    "importRecords(var,sourceDb,sql,table)",
    "purge(var)",
    "runSQLs(var,sqls)"
)
for ($i = 0; $i -lt $localdb.Length; $i++) {
    $ws.Cells.Item(1 + $i, 14).Value = $localdb[$i]
}

# 3) Insert "localdb" into the alphabetically-sorted category list in column
#    A (the "target" named range), right before "macro", pushing the
#    remaining entries (macro, mail, number, pdf, rdbms, redis, sms, sound,
#    ssh, step, web, webalert, webcookie, ws, ws.async, xml) down by one row.
$target = @(
    "localdb",
    "macro",
    "mail",
    "number",
    "pdf",
    "rdbms",
    "redis",
    "sms",
    "sound",
    "ssh",
    "step",
    "web",
    "webalert",
    "webcookie",
    "ws",
    "ws.async",
    "xml"
)
for ($i = 0; $i -lt $target.Length; $i++) {
    $ws.Cells.Item(14 + $i, 1).Value = $target[$i]
}

# 4) Fix up the named ranges that pointed at the now-shifted columns, and add
#    the brand-new "localdb" name.
$wb.Names.Item("mail").RefersTo = "='#system'!`$P`$2:`$P`$2"
$wb.Names.Item("number").RefersTo = "='#system'!`$Q`$2:`$Q`$16"
$wb.Names.Item("pdf").RefersTo = "='#system'!`$R`$2:`$R`$16"
$wb.Names.Item("rdbms").RefersTo = "='#system'!`$S`$2:`$S`$7"
$wb.Names.Item("redis").RefersTo = "='#system'!`$T`$2:`$T`$10"
$wb.Names.Item("sms").RefersTo = "='#system'!`$U`$2:`$U`$2"
$wb.Names.Item("sound").RefersTo = "='#system'!`$V`$2:`$V`$5"
$wb.Names.Item("ssh").RefersTo = "='#system'!`$W`$2:`$W`$9"
$wb.Names.Item("step").RefersTo = "='#system'!`$X`$2:`$X`$4"
$wb.Names.Item("web").RefersTo = "='#system'!`$Y`$2:`$Y`$127"
$wb.Names.Item("webalert").RefersTo = "='#system'!`$Z`$2:`$Z`$8"
$wb.Names.Item("webcookie").RefersTo = "='#system'!`$AA`$2:`$AA`$8"
$wb.Names.Item("ws").RefersTo = "='#system'!`$AB`$2:`$AB`$17"
$wb.Names.Item("ws.async").RefersTo = "='#system'!`$AC`$2:`$AC`$8"
$wb.Names.Item("xml").RefersTo = "='#system'!`$AD`$2:`$AD`$21"
$wb.Names.Item("macro").RefersTo = "='#system'!`$O`$2:`$O`$4"
$wb.Names.Item("target").RefersTo = "='#system'!`$A`$2:`$A`$30"
$wb.Names.Add("localdb", "='#system'!`$N`$2:`$N`$7")
